# Update the two-digit-divided-by-one-digit practice problems in the
# single table on the page.
#
# Several of the original problem strings are duplicated across cells
# (e.g. "71÷8=" occurs twice, but must become two different results),
# so a document-wide Find/Replace would be ambiguous. Each cell is
# therefore addressed directly by its table row/column, and only the
# visible text portion of the cell (excluding the trailing paragraph
# mark and cell-end mark) is replaced.

function Set-CellText {
    param($Table, $Row, $Col, $OldText, $NewText)
    $cell = $Table.Cell($Row, $Col)
    $rng = $cell.Range
    $rng.End = $rng.Start + $OldText.Length
    $rng.Text = $NewText
}

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

Set-CellText $t 1 1 "33÷5=" "34÷4="
Set-CellText $t 1 2 "71÷8=" "69÷7="
Set-CellText $t 1 3 "23÷4=" "53÷2="
Set-CellText $t 1 4 "34÷6=" "67÷3="
Set-CellText $t 1 5 "61÷3=" "48÷5="

Set-CellText $t 5 1 "50÷8=" "86÷2="
Set-CellText $t 5 2 "71÷8=" "44÷2="
Set-CellText $t 5 3 "29÷5=" "10÷7="
Set-CellText $t 5 4 "13÷8=" "62÷2="
Set-CellText $t 5 5 "10÷6=" "67÷7="

Set-CellText $t 9 1 "87÷6=" "25÷8="
Set-CellText $t 9 2 "21÷6=" "29÷5="
Set-CellText $t 9 3 "47÷8=" "21÷5="
Set-CellText $t 9 4 "94÷6=" "51÷7="
Set-CellText $t 9 5 "63÷2=" "31÷6="

Set-CellText $t 13 1 "76÷7=" "69÷6="
Set-CellText $t 13 2 "97÷3=" "63÷3="
Set-CellText $t 13 3 "57÷3=" "76÷5="
Set-CellText $t 13 4 "88÷7=" "50÷6="
Set-CellText $t 13 5 "66÷6=" "56÷7="

Set-CellText $t 17 1 "37÷8=" "21÷4="
Set-CellText $t 17 2 "43÷6=" "25÷8="
Set-CellText $t 17 3 "54÷8=" "97÷4="
Set-CellText $t 17 4 "96÷4=" "32÷5="
Set-CellText $t 17 5 "14÷6=" "38÷6="

Write-Host "Done updating problems."
